$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

$ws.Range('D2').Value = '63.637.19'
$ws.Range('E2').Value = '  -0.93%  '
$ws.Range('D3').Value = '2.612.81'
$ws.Range('E3').Value = '  -0.50%  '
Set-TextValue 'D5' '590.95'
$ws.Range('E5').Value = '  -1.85%  '
Set-TextValue 'D6' '150.07'
$ws.Range('E6').Value = '  -0.42%  '
$ws.Range('E7').Value = '  -0.01%  '
Set-TextValue 'D8' '0.585'
$ws.Range('E8').Value = '  -0.86%  '
$ws.Range('E9').Value = '  +0.01%  '
Set-TextValue 'D10' '5.76'
$ws.Range('E10').Value = '  +0.82%  '
$ws.Range('E11').Value = '  +0.35%  '
$ws.Range('E12').Value = '  +0.38%  '
Set-TextValue 'D13' '27.71'
$ws.Range('E13').Value = '  -0.01%  '
$ws.Range('D14').Value = '3.080.75'
$ws.Range('E14').Value = '  -0.57%  '
$ws.Range('D15').Value = '63.459.04'
$ws.Range('E15').Value = '  -0.96%  '
$ws.Range('E16').Value = '  +4.90%  '
$ws.Range('D17').Value = '2.608.93'
$ws.Range('E17').Value = '  -0.83%  '
Set-TextValue 'D18' '12.17'
$ws.Range('E18').Value = '  -0.17%  '
Set-TextValue 'D19' '4.77'
$ws.Range('E19').Value = '  +2.31%  '
Set-TextValue 'D20' '345.86'
$ws.Range('E20').Value = '  -1.62%  '
Set-TextValue 'D21' '6.90'
$ws.Range('E21').Value = '  -1.07%  '
$ws.Range('E22').Value = '  +0.13%  '
Set-TextValue 'D23' '66.99'
$ws.Range('E23').Value = '  +0.53%  '
$ws.Range('E24').Value = '  -3.00%  '
Set-TextValue 'D25' '9.24'
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('E26').Value = '  -1.92%  '
Set-TextValue 'D27' '8.46'
$ws.Range('E27').Value = '  +3.51%  '
Set-TextValue 'D28' '548.83'
$ws.Range('E28').Value = '  +2.02%  '
$ws.Range('E29').Value = '  -2.75%  '
$ws.Range('E30').Value = '  -0.05%  '
Set-TextValue 'D31' '2.04'
$ws.Range('E31').Value = '  -1.17%  '
$ws.Range('D32').Value = '0.0₃0870'
$ws.Range('E32').Value = '  +1.69%  '
$ws.Range('E33').Value = '  +1.53%  '
Set-TextValue 'D34' '5.36'
$ws.Range('E34').Value = '  +1.25%  '
Set-TextValue 'D35' '6.13'
$ws.Range('E35').Value = '  -0.69%  '
Set-TextValue 'D36' '165.92'
$ws.Range('E36').Value = '  -1.21%  '
Set-TextValue 'D37' '0.413'
$ws.Range('E37').Value = '  +0.70%  '
$ws.Range('E38').Value = '  -0.14%  '
$ws.Range('E39').Value = '  -1.89%  '
Set-TextValue 'D40' '19.54'
$ws.Range('E40').Value = '  +0.11%  '
$ws.Range('E41').Value = '  -0.02%  '
Set-TextValue 'D42' '165.31'
$ws.Range('E42').Value = '  -2.62%  '
Set-TextValue 'D43' '4.05'
$ws.Range('E43').Value = '  +2.69%  '
Set-TextValue 'D44' '23.15'
$ws.Range('E44').Value = '  +7.14%  '
Set-TextValue 'D45' '0.0582'
$ws.Range('E45').Value = '  -2.03%  '
$ws.Range('E46').Value = '  +7.73%  '
Set-TextValue 'D47' '0.633'
$ws.Range('E47').Value = '  +0.56%  '
$ws.Range('E48').Value = '  +1.88%  '
Set-TextValue 'D49' '0.0962'
$ws.Range('E49').Value = '  -0.69%  '
$ws.Range('E50').Value = '  -0.52%  '
$ws.Range('D51').Value = '0.0₆0233'
$ws.Range('E51').Value = '  +17.27%  '
